$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the two newly-added checklist items as accomplished ("X")
$ws.Range("B15").Value = "X"
$ws.Range("B16").Value = "X"

# Update the active selection to B5
$ws.Range("B5").Select()
